$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labels (rows 29-32, column A) -- shared strings 20..23
$ws.Range("A29").Value = "Période de récupération"
$ws.Range("A30").Value = "Avantages-Investissement"
$ws.Range("A31").Value = "Investissement"
$ws.Range("A32").Value = "Rendement du capital investi sur 5 ans"

# New values / formulas (column B)
$ws.Range("B29").Value = 1
$ws.Range("B30").Formula = "=G27"
$ws.Range("B31").Formula = "=-(B24+C25+D25+E25+F25+G25)"
$ws.Range("B32").Formula = "=B30/B31"

# B30/B31 show the investment totals in the same currency ("Monétaire") look
# used throughout column B; B32 is the 5-year ROI expressed as a percentage.
$ws.Range("B30:B31").NumberFormat = "_ * #,##0.00_)\ ""$""_ ;_ * \(#,##0.00\)\ ""$""_ ;_ * ""-""??_)\ ""$""_ ;_ @_ "
$ws.Range("B32").NumberFormat = "0.00%"

# Move the selection to the newly-added summary cell, like in the saved file
$ws.Range("B32").Select()
